$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 105, shifting existing rows 105:110 down to 106:111.
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row 105 with the new weekly price record.
$ws.Cells.Item(105, 1).Value = 1
$ws.Cells.Item(105, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(105, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(105, 4).Value = 44858
$ws.Cells.Item(105, 4).NumberFormat = $ws.Cells.Item(106, 4).NumberFormat
$ws.Cells.Item(105, 5).Value = 15
$ws.Cells.Item(105, 6).Value = 100112021
$ws.Cells.Item(105, 7).Value = "Ají"
$ws.Cells.Item(105, 8).Value = "Cristal"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 250
$ws.Cells.Item(105, 11).Value = 24000
$ws.Cells.Item(105, 12).Value = 25000
$ws.Cells.Item(105, 13).Value = 24400
$ws.Cells.Item(105, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(105, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(105, 16).Value = 1627
$ws.Cells.Item(105, 17).Value = 15
$ws.Cells.Item(105, 18).Value = "Hortaliza"
